$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 12.83479999999999
$ws.Range("E6").Value = 12.4856
$ws.Range("D7").Value = -7.324199999999993
$ws.Range("A10").Value = -20.53679999999997
$ws.Range("A12").Value = -22.71160000000004
$ws.Range("B13").Value = 6.008799999999997
$ws.Range("A18").Value = -22.60950000000003
$ws.Range("D20").Value = -8.456899999999997
